$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - Excel recalculated / re-saved the timestamp with a (sub-millisecond)
# different floating point value for column A.
$ws.Range("A14").Value = 45878.58356094908

# New row 15 with the latest weather reading.
$ws.Range("A15").Value = 45878.62517811374
$ws.Range("A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B15").Value = 2025
$ws.Range("C15").Value = 37
$ws.Range("D15").Value = 20.11
$ws.Range("E15").Value = 74.53
$ws.Range("F15").Value = 446.48
$ws.Range("G15").Value = 14.94
$ws.Range("H15").Value = "ESE"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "15:00:15"
